{"js": "// Apply tense/phrasing changes to the abstract body paragraph.\n// We perform a sequence of targeted search-and-replace operations that,\n// together, transform the original text into the revised text described\n// by the diff (changing verb tense and rewording a few sentences near\n// the end of the paragraph).\n\nconst body = context.document.body;\n\nasync function replaceOnce(searchText, replacement) {\n  const results = body.search(searchText, { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Search text not found: \" + searchText);\n  }\n\n  results.items[0].insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1) \"We use mathematical modeling\" -> \"We used mathematical modeling\"\nawait replaceOnce(\n  \"We use mathematical modeling\",\n  \"We used mathematical modeling\"\n);\n\n// 2) \"A family of GRNs has been derived\" -> \"A family of GRNs was derived\"\nawait replaceOnce(\n  \"A family of GRNs has been derived\",\n  \"A family of GRNs was derived\"\n);\n\n// 3) \"we will compare the modeling results\" -> \"we compared the modeling results\"\nawait replaceOnce(\n  \"we will compare the modeling results\",\n  \"we compared the modeling results\"\n);\n\n// 4) Rework of the final few sentences (future tense -> past tense, plus\n// rewording of the comparison/analysis/validation sentences). This is\n// split into two replacements around the \"_GoBack\" bookmark (which sits\n// in the middle of \"DB-derived\") so that the bookmark is preserved in\n// its original location rather than being collapsed/moved.\nawait replaceOnce(\n  \"a large collection of random networks will be generated via an R script. Compariso\",\n  \"a large collection of random networks was generated. Comparisons made between the random networks and the DB-derived network consistently saw better modeling of the DB-der\"\n);\n\nawait replaceOnce(\n  \"ns will be made between the random networks and the DB-derived network. \" +\n    \"We predict that we will see a significant difference between the random network and the DB-derived network. \" +\n    \"We will analyze the network to determine which features better predict yeast cell behavior. \" +\n    \"This in turn will validate our predictions\",\n  \"ived network. Through analysis, we determined key features of the DB-derived network that better predicted yeast cell behavior. \" +\n    \"This in turn validated our predictions\"\n);\n", "ps1": "# Apply tense/phrasing changes to the abstract body paragraph.\n# We perform a sequence of targeted Find & Replace operations that,\n# together, transform the original text into the revised text described\n# by the diff (changing verb tense and rewording a few sentences near\n# the end of the paragraph).\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $found = $range.Find.Execute(\n        $findText,    # FindText\n        $false,       # MatchCase\n        $false,       # MatchWholeWord\n        $false,       # MatchWildcards\n        $false,       # MatchSoundsLike\n        $false,       # MatchAllWordForms\n        $true,        # Forward\n        1,            # Wrap (wdFindContinue)\n        $false,       # Format\n        $replaceText, # ReplaceWith\n        2             # Replace (wdReplaceOne)\n    )\n    if (-not $found) {\n        throw \"Could not find text to replace: $findText\"\n    }\n}\n\n# 1) \"We use mathematical modeling\" -> \"We used mathematical modeling\"\nReplace-Text \"We use mathematical modeling\" \"We used mathematical modeling\"\n\n# 2) \"A family of GRNs has been derived\" -> \"A family of GRNs was derived\"\nReplace-Text \"A family of GRNs has been derived\" \"A family of GRNs was derived\"\n\n# 3) \"we will compare the modeling results\" -> \"we compared the modeling results\"\nReplace-Text \"we will compare the modeling results\" \"we compared the modeling results\"\n\n# 4) Rework of the final few sentences (future tense -> past tense, plus\n# rewording of the comparison/analysis/validation sentences). This is\n# split into two replacements around the \"_GoBack\" bookmark (which sits\n# in the middle of \"DB-derived\") so that the bookmark is preserved in its\n# original location rather than being collapsed/removed by the replace.\nReplace-Text \"a large collection of random networks will be generated via an R script. Compariso\" \"a large collection of random networks was generated. Comparisons made between the random networks and the DB-derived network consistently saw better modeling of the DB-der\"\n\nReplace-Text \"ns will be made between the random networks and the DB-derived network. We predict that we will see a significant difference between the random network and the DB-derived network. We will analyze the network to determine which features better predict yeast cell behavior. This in turn will validate our predictions\" \"ived network. Through analysis, we determined key features of the DB-derived network that better predicted yeast cell behavior. This in turn validated our predictions\"\n"}
